# Applies the betexplorer data-refresh edit:
#  1. Eight pairs of adjacent rows have their match data (columns B:V)
#     swapped between the two rows (the "Indice" in column A stays put,
#     tied to the physical row).
#  2. Two brand-new match rows (112 and 113) are appended at the bottom,
#     extending the used range from A1:V111 to A1:V113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($r1, $r2) {
    $rangeA = $ws.Range("B$r1`:V$r1")
    $rangeB = $ws.Range("B$r2`:V$r2")
    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

# Row pairs whose match data needs to be exchanged.
Swap-MatchRows 2 3
Swap-MatchRows 4 5
Swap-MatchRows 26 27
Swap-MatchRows 28 29
Swap-MatchRows 53 54
Swap-MatchRows 57 58
Swap-MatchRows 59 60
Swap-MatchRows 86 87

# Append the two new rows (112, 113) at the end of the sheet, copying the
# formatting from the last existing row (111) so the index/date columns
# keep their number formats.
$ws.Range("A111:V111").Copy()
$ws.Range("A112:V113").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

function Set-MatchRow($r, $vals) {
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

Set-MatchRow 112 @(111, "italy", "serie-a", "2023-2024", 45240.77083333334, "Sassuolo", 2, "Salernitana", 2, 1.76, "29/10/2023 11:02", 1.76, "10/11/2023 18:29", 3.97, "29/10/2023 11:02", 4, "10/11/2023 18:29", 4.52, "29/10/2023 11:02", 4.7, "10/11/2023 18:29", "https://www.betexplorer.com/football/italy/serie-a/sassuolo-salernitana/z9JgdBup/")

Set-MatchRow 113 @(112, "italy", "serie-a", "2023-2024", 45240.86458333334, "Genoa", 1, "Verona", 0, 1.95, "29/10/2023 11:02", 1.87, "10/11/2023 20:44", 3.43, "29/10/2023 11:02", 3.35, "10/11/2023 20:44", 4.21, "29/10/2023 11:02", 5.09, "10/11/2023 20:44", "https://www.betexplorer.com/football/italy/serie-a/genoa-verona/z9s26e3i/")
